# Kanban board update: "Sign-up System" and its preceding setup tasks moved
# from "Not Started" (column A) to "Done" (column C); "Create cloud
# infrastructure" moved from "Doing" (column B) to "Done" (column C) too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("Not Started") -------------------------------------------
# Rows 4-8 (Initialize Backend Folder, Add entities in MySQL, Initialize
# Frontend Folder, Initialize AWS Environment, Sign-up System) are done, so
# everything below shifts up by 5 rows; the newly-vacated tail (15-19)
# becomes blank.
$ws.Range("A4").Value  = "Login System"
$ws.Range("A5").Value  = "Password Recovery System"
$ws.Range("A6").Value  = "Homepage System"
$ws.Range("A7").Value  = "Account System: Edit Profile"
$ws.Range("A8").Value  = "Account System: Change Password"
$ws.Range("A9").Value  = "Account System: Change Email"
$ws.Range("A10").Value = "Add Books System"
$ws.Range("A11").Value = "Update Books System"
$ws.Range("A12").Value = "User Idle System"
$ws.Range("A13").Value = "Website UI"
$ws.Range("A14").Value = "View Book Instance Page"
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()
$ws.Range("A17").ClearContents()
$ws.Range("A18").ClearContents()
$ws.Range("A19").ClearContents()

# --- Column B ("Doing") ---------------------------------------------------
# "Create cloud infrastructure" is finished, so it leaves the Doing column.
$ws.Range("B4").ClearContents()

# --- Column C ("Done") -----------------------------------------------------
# Newly finished items are appended after the existing two (rows 4 & 5 stay
# the same): Create cloud infrastructure, then the four setup tasks, then
# Sign-up System.
$ws.Range("C6").Value  = "Create cloud infrastructure"
$ws.Range("C7").Value  = "Initialize Backend Folder"
$ws.Range("C8").Value  = "Add entities in MySQL"
$ws.Range("C9").Value  = "Initialize Frontend Folder"
$ws.Range("C10").Value = "Initialize AWS Environment"
$ws.Range("C11").Value = "Sign-up System"

# --- A3 formula tweak --------------------------------------------------
# The COUNTA range start for "Not Started" moved from A10 to A5 (still
# counts every populated row beneath the header).
$ws.Range("A3").Formula = "=COUNTA(A5:A1048576)"

# --- View state: selection + zoom ------------------------------------------
$ws.Range("B13").Select()
$excel.ActiveWindow.Zoom = 85
